$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '24.126.44'
    'E2' = '  -3.03%  '
    'D3' = '1.641.09'
    'E3' = '  -2.89%  '
    'E4' = '  -0.37%  '
    'D5' = '306.54'
    'E5' = '  -2.55%  '
    'D6' = '1.000'
    'E6' = '  -0.37%  '
    'D7' = '0.3889'
    'E7' = '  -1.23%  '
    'D8' = '0.3861'
    'E8' = '  -3.21%  '
    'D9' = '0.9994'
    'E9' = '  -0.52%  '
    'D10' = '49.47'
    'E10' = '  -5.51%  '
    'D11' = '1.349'
    'E11' = '  -6.18%  '
    'D12' = '0.08702'
    'E12' = '  -0.11%  '
    'D13' = '23.64'
    'E13' = '  -7.09%  '
    'D14' = '7.100'
    'E14' = '  -3.45%  '
    'D15' = '0.00001288'
    'E15' = '  -3.33%  '
    'D16' = '7.457'
    'E16' = '  -4.97%  '
    'D17' = '1.632.99'
    'E17' = '  +2.02%  '
    'D18' = '95.14'
    'E18' = '  +0.61%  '
    'D19' = '0.06892'
    'E19' = '  -3.82%  '
    'D20' = '20.58'
    'E20' = '  +1.06%  '
    'D21' = '6.896'
    'E21' = '  -4.01%  '
    'D22' = '0.9992'
    'E22' = '  -0.44%  '
    'D23' = '13.55'
    'E23' = '  -4.43%  '
    'D24' = '24.140.50'
    'E24' = '  -2.97%  '
    'D25' = '2.327'
    'E25' = '  -2.66%  '
    'D26' = '2.736'
    'E26' = '  -3.78%  '
    'D27' = '22.32'
    'E27' = '  -2.74%  '
    'D28' = '157.57'
    'E28' = '  -2.61%  '
    'D29' = '8.531'
    'E29' = '  +6.45%  '
    'D30' = '140.02'
    'E30' = '  -5.40%  '
    'D31' = '5.344'
    'E31' = '  -10.47%  '
    'D32' = '2.411'
    'E32' = '  -8.37%  '
    'D33' = '1.816.55'
    'E33' = '  +0.12%  '
    'D34' = '6.916'
    'E34' = '  -1.17%  '
    'D35' = '0.07998'
    'E35' = '  -6.11%  '
    'D36' = '0.02879'
    'E36' = '  -7.01%  '
    'D37' = '0.2671'
    'E37' = '  -6.65%  '
    'E38' = '  -7.88%  '
    'D39' = '0.09180'
    'E39' = '  -4.85%  '
    'D40' = '1.467'
    'E40' = '  +0.28%  '
    'D41' = '9.908'
    'E41' = '  -7.52%  '
    'D42' = '0.7543'
    'E42' = '  -6.34%  '
    'D43' = '13.02'
    'E43' = '  -5.86%  '
    'E44' = '  -4.82%  '
    'D45' = '0.6884'
    'E45' = '  -5.00%  '
    'D46' = '2.465'
    'E46' = '  -5.95%  '
    'D47' = '4.085'
    'E47' = '  -3.04%  '
    'D48' = '0.9994'
    'E48' = '  -0.26%  '
    'D49' = '0.08387'
    'E49' = '  -5.89%  '
    'D50' = '1.259'
    'E50' = '  -8.97%  '
    'D51' = '132.87'
    'E51' = '  -4.34%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
